# Re-run SGNN to annotate dialog acts following clean up work to the original transcripts.
# Update DAMSLTag (column I) and DialogAct (column J) values for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 5;  I = "sv"; J = "Statement-opinion" },
    @{ Row = 10; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 14; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 19; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 33; I = "aa"; J = "Agree/Accept" },
    @{ Row = 42; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 43; I = "sv"; J = "Statement-opinion" },
    @{ Row = 44; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 57; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 60; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 62; I = "ba"; J = "Appreciation" },
    @{ Row = 64; I = "aa"; J = "Agree/Accept" },
    @{ Row = 65; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 67; I = "ba"; J = "Appreciation" }
)

foreach ($u in $updates) {
    $ws.Range("I" + $u.Row).Value = $u.I
    $ws.Range("J" + $u.Row).Value = $u.J
}
